$wb = $excel.ActiveWorkbook

# --- Yearly sheet: November's Taxable Account dividend was recorded (was 0) ---
$wsYearly = $wb.Worksheets.Item("Yearly")
$wsYearly.Range("D13").Value = 66.15
# (G13, D15, G15 totals on this sheet recalc automatically from the formulas already in place)

# Leave the cursor where the author left it when they saved
$wsYearly.Activate()
$null = $wsYearly.Range("G18").Select()

# --- All Time sheet: just the view/selection changed; values flow in via formulas ---
$wsAllTime = $wb.Worksheets.Item("All Time")
$wsAllTime.Activate()
$null = $wsAllTime.Range("H14").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
